# Reshape the sheet to match a pandas re-export that carries an extra
# "Unnamed: 0" index column in front of the previous data.
#
# Old layout: A=id, B=processingdate, C=transcription, D=summary_old,
#             E=summary, F=topic, G=sentiment, H=cost
# New layout: A=(blank header, 0-based index), B=Unnamed: 0 (0-based index),
#             C=id, D=processingdate, E=transcription, F=summary_old,
#             G=summary, H=topic, I=sentiment, J=cost

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A; this shifts old A:H -> C:J
# (content, types and formatting all travel with the shift).
$ws.Range("A:B").Insert()

# New header cells get the same bold/centered/bordered style used by the
# rest of the header row (copy format from the "id" header, now in C1).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Column B header text.
$ws.Range("B1").Value = 'Unnamed: 0'

# Columns A and B both hold the new 0-based row index.
for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $i
}

# The old "id" column (now column C) becomes numeric instead of text.
$idValues = @(18531, 18526, 18522, 18509, 18498, 18497, 18496, 18492, 18491, 18467)
for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $idValues[$i]
}
